$d = $word.ActiveDocument

function Replace-Text($old, $new) {
  $ok = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
  if (-not $ok) {
    Write-Output ("WARNING: replace failed for: " + $old)
  }
  return $ok
}

# Stack-trace line-number bumps (M2DocEvaluator / AbstractTemplatesTestSuite)
# caused by the M2Doc source moving from 2.0.1 to 2.0.2.
Replace-Text "M2DocEvaluator.java:555)" "M2DocEvaluator.java:559)"
Replace-Text "M2DocEvaluator.java:1096)" "M2DocEvaluator.java:1216)"
Replace-Text "M2DocEvaluator.java:1305)" "M2DocEvaluator.java:1425)"
Replace-Text "M2DocEvaluator.java:283)" "M2DocEvaluator.java:287)"
Replace-Text "M2DocEvaluator.java:272)" "M2DocEvaluator.java:276)"
Replace-Text "AbstractTemplatesTestSuite.java:479)" "AbstractTemplatesTestSuite.java:480)"
Replace-Text "AbstractTemplatesTestSuite.java:388)" "AbstractTemplatesTestSuite.java:389)"
Replace-Text "GeneratedMethodAccessor75" "GeneratedMethodAccessor74"

$oldTail = "	at org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)" + [char]10 + `
"	at org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)" + [char]10 + `
"	at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:539)" + [char]10 + `
"	at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:761)" + [char]10 + `
"	at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:461)" + [char]10 + `
"	at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:207)"

$newTail = "	at org.apache.maven.surefire.junit4.JUnit4Provider.execute(JUnit4Provider.java:264)" + [char]10 + `
"	at org.apache.maven.surefire.junit4.JUnit4Provider.executeTestSet(JUnit4Provider.java:153)" + [char]10 + `
"	at org.apache.maven.surefire.junit4.JUnit4Provider.invoke(JUnit4Provider.java:124)" + [char]10 + `
"	at sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)" + [char]10 + `
"	at sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)" + [char]10 + `
"	at sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)" + [char]10 + `
"	at java.lang.reflect.Method.invoke(Method.java:498)" + [char]10 + `
"	at org.apache.maven.surefire.util.ReflectionUtils.invokeMethodWithArray2(ReflectionUtils.java:208)" + [char]10 + `
"	at org.apache.maven.surefire.booter.ProviderFactory`$ProviderProxy.invoke(ProviderFactory.java:156)" + [char]10 + `
"	at org.apache.maven.surefire.booter.ProviderFactory.invokeProvider(ProviderFactory.java:82)" + [char]10 + `
"	at org.eclipse.tycho.surefire.osgibooter.OsgiSurefireBooter.run(OsgiSurefireBooter.java:91)" + [char]10 + `
"	at org.eclipse.tycho.surefire.osgibooter.HeadlessTestApplication.run(HeadlessTestApplication.java:21)" + [char]10 + `
"	at sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)" + [char]10 + `
"	at sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)" + [char]10 + `
"	at sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)" + [char]10 + `
"	at java.lang.reflect.Method.invoke(Method.java:498)" + [char]10 + `
"	at org.eclipse.equinox.internal.app.EclipseAppContainer.callMethodWithException(EclipseAppContainer.java:587)" + [char]10 + `
"	at org.eclipse.equinox.internal.app.EclipseAppHandle.run(EclipseAppHandle.java:198)" + [char]10 + `
"	at org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.runApplication(EclipseAppLauncher.java:134)" + [char]10 + `
"	at org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.start(EclipseAppLauncher.java:104)" + [char]10 + `
"	at org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:388)" + [char]10 + `
"	at org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:243)" + [char]10 + `
"	at sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)" + [char]10 + `
"	at sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)" + [char]10 + `
"	at sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)" + [char]10 + `
"	at java.lang.reflect.Method.invoke(Method.java:498)" + [char]10 + `
"	at org.eclipse.equinox.launcher.Main.invokeFramework(Main.java:656)" + [char]10 + `
"	at org.eclipse.equinox.launcher.Main.basicRun(Main.java:592)" + [char]10 + `
"	at org.eclipse.equinox.launcher.Main.run(Main.java:1498)" + [char]10 + `
"	at org.eclipse.equinox.launcher.Main.main(Main.java:1471)"


# Replace the Eclipse JDT test-runner tail of the stack trace with the
# Maven Surefire / Tycho / Equinox launcher tail (new CI runner).
Replace-Text $oldTail $newTail
